# Update the vendor names in column A for the three "personal" rows,
# renaming them to match the new naming scheme used elsewhere in the sheet.
# The corresponding addresses in column B are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Eupeptic"   # was "Dad"     -> 349 Eupeptic Springs Rd, North Carolina
$ws.Range("A14").Value = "Cabin"      # was "Mamaw"   -> 126 Sharpe Bluff Ln, North Carolina
$ws.Range("A15").Value = "Apartments" # was "Kinsley" -> 121 Pine Valley Dr, Yadkinville, NC 27055
